$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 1-15, 17-31, 33 (values only; B/C already Text-formatted) ---
# Row 1
$ws.Range("A1").Value2 = "Hare"
$ws.Range("B1").Value2 = "00"
$ws.Range("C1").Value2 = "175A89F4"
$ws.Range("G1").Value2 = "07292AF4"

# Row 2
$ws.Range("A2").Value2 = "Monkey"
$ws.Range("B2").Value2 = "01"
$ws.Range("C2").Value2 = "F7C68DF4"

# Row 3
$ws.Range("A3").Value2 = "Penguin"
$ws.Range("B3").Value2 = "02"
$ws.Range("C3").Value2 = "D7965CF4"

# Row 4
$ws.Range("A4").Value2 = "Hippopotamus"
$ws.Range("B4").Value2 = "03"
$ws.Range("C4").Value2 = "A49E502E"

# Row 5
$ws.Range("A5").Value2 = "Sheep"
$ws.Range("B5").Value2 = "04"
$ws.Range("C5").Value2 = "87DC8CF4"

# Row 6
$ws.Range("A6").Value2 = "Rabbit"
$ws.Range("B6").Value2 = "05"
$ws.Range("C6").Value2 = "74A8722E"

# Row 7
$ws.Range("A7").Value2 = "Dog"
$ws.Range("B7").Value2 = "06"
$ws.Range("C7").Value2 = "F421682E"

# Row 8
$ws.Range("A8").Value2 = "Frog"
$ws.Range("B8").Value2 = "07"
$ws.Range("C8").Value2 = "674132F4"

# Row 9
$ws.Range("A9").Value2 = "Bee"
$ws.Range("B9").Value2 = "08"
$ws.Range("C9").Value2 = "87AA36F4"
$ws.Range("G9").Value2 = "87AA36F4"

# Row 10
$ws.Range("A10").Value2 = "Elephant"
$ws.Range("B10").Value2 = "09"
$ws.Range("C10").Value2 = "54D06A2E"
$ws.Range("G10").Value2 = "54D06A2E"

# Row 11
$ws.Range("A11").Value2 = "GingerBreadMan"
$ws.Range("B11").Value2 = "10"
$ws.Range("C11").Value2 = "A7B243F4"

# Row 12
$ws.Range("A12").Value2 = "Chicken"
$ws.Range("B12").Value2 = "11"
$ws.Range("C12").Value2 = "978B35F4"

# Row 13
$ws.Range("A13").Value2 = "Chick"
$ws.Range("B13").Value2 = "12"
$ws.Range("C13").Value2 = "74CB662E"

# Row 14
$ws.Range("A14").Value2 = "Hedgehog"
$ws.Range("B14").Value2 = "13"
$ws.Range("C14").Value2 = "B7538CF4"

# Row 15
$ws.Range("A15").Value2 = "Mouse"
$ws.Range("B15").Value2 = "14"
$ws.Range("C15").Value2 = "773536F4"

# Row 17
$ws.Range("A17").Value2 = "Hare"
$ws.Range("B17").Value2 = "00"
$ws.Range("C17").Value2 = "07292AF4"

# Row 18
$ws.Range("A18").Value2 = "Monkey"
$ws.Range("B18").Value2 = "01"
$ws.Range("C18").Value2 = "775228F4"

# Row 19
$ws.Range("A19").Value2 = "Penguin"
$ws.Range("B19").Value2 = "02"
$ws.Range("C19").Value2 = "671931F4"

# Row 20
$ws.Range("A20").Value2 = "Hippopotamus"
$ws.Range("B20").Value2 = "03"
$ws.Range("C20").Value2 = "97602CF4"

# Row 21
$ws.Range("A21").Value2 = "Sheep"
$ws.Range("B21").Value2 = "04"
$ws.Range("C21").Value2 = "A7748EF4"

# Row 22
$ws.Range("A22").Value2 = "Rabbit"
$ws.Range("B22").Value2 = "05"
$ws.Range("C22").Value2 = "E77090F4"

# Row 23
$ws.Range("A23").Value2 = "Dog"
$ws.Range("B23").Value2 = "06"
$ws.Range("C23").Value2 = "07A635F4"

# Row 24
$ws.Range("A24").Value2 = "Frog"
$ws.Range("B24").Value2 = "07"
$ws.Range("C24").Value2 = "675E90F4"

# Row 25
$ws.Range("A25").Value2 = "Bee"
$ws.Range("B25").Value2 = "08"
$ws.Range("C25").Value2 = "27BF63F4"

# Row 26
$ws.Range("A26").Value2 = "Elephant"
$ws.Range("B26").Value2 = "09"
$ws.Range("C26").Value2 = "371E83F4"

# Row 27
$ws.Range("A27").Value2 = "GingerBreadMan"
$ws.Range("B27").Value2 = "10"
$ws.Range("C27").Value2 = "37642BF4"

# Row 28
$ws.Range("A28").Value2 = "Chicken"
$ws.Range("B28").Value2 = "11"
$ws.Range("C28").Value2 = "D4DD552E"

# Row 29
$ws.Range("A29").Value2 = "Chick"
$ws.Range("B29").Value2 = "12"
$ws.Range("C29").Value2 = "E7145FF4"

# Row 30
$ws.Range("A30").Value2 = "Hedgehog"
$ws.Range("B30").Value2 = "13"
$ws.Range("C30").Value2 = "E4404F2E"

# Row 31
$ws.Range("A31").Value2 = "Mouse"
$ws.Range("B31").Value2 = "14"
$ws.Range("C31").Value2 = "94785B2E"

# Row 33
$ws.Range("A33").Value2 = "Reprint"
$ws.Range("B33").Value2 = "100"
$ws.Range("C33").Value2 = "84AA732E"

# --- Add new rows 34-37 (admin commands) ---
# Row 34
$ws.Range("A34").Value2 = "AbortGame"
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value2 = "101"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value2 = "67668FF4"
$ws.Range("D34").Formula = '=CONCATENATE(C34,CHAR(9),B34,"|",A34)'
$ws.Range("E34").Formula = '=LEFT(A34,16)'

# Row 35
$ws.Range("A35").Value2 = "Player"
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value2 = "102"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value2 = "677C36F4"
$ws.Range("D35").Formula = '=CONCATENATE(C35,CHAR(9),B35,"|",A35)'
$ws.Range("E35").Formula = '=LEFT(A35,16)'

# Row 36
$ws.Range("A36").Value2 = "Rounds"
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value2 = "103"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value2 = "4460792E"
$ws.Range("D36").Formula = '=CONCATENATE(C36,CHAR(9),B36,"|",A36)'
$ws.Range("E36").Formula = '=LEFT(A36,16)'

# Row 37
$ws.Range("A37").Value2 = "Show"
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value2 = "104"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value2 = "77B322F4"
$ws.Range("D37").Formula = '=CONCATENATE(C37,CHAR(9),B37,"|",A37)'
$ws.Range("E37").Formula = '=LEFT(A37,16)'

# --- View state: selection + best-effort scroll position ---
$ws.Range("C33").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1

Write-Output "edit complete"
